$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pcb-adr1399-vref-kicost")

$ws.Range("B3").Value = "Fri Feb  4 23:08:02 2022"
$ws.Range("B4").Value = "2022-02-04 23:08:03"

$ws.Range("C24").Value = 1.182130909176882
$ws.Range("C25").Value = 0.8722958827634333
